$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# The "Gruppe" column (B) duplicated the group label that already lives in
# the last column, so select it (as a user would before deleting it) and
# remove it, shifting everything after it one column to the left.
[void]$ws1.Columns.Item(2).Select()
[void]$ws1.Columns.Item(2).Delete()

# Re-point the leftover (stale) manual-sort bookkeeping at the columns it
# referred to before the deletion (D -> C, R -> Q), without re-sorting the
# two sample rows that are actually present. We briefly equalize the sort
# key so Apply() can't reorder anything, then restore the real values.
$origC4 = $ws1.Range("C4").Value()
$origC5 = $ws1.Range("C5").Value()
$ws1.Range("C4").Value = 1
$ws1.Range("C5").Value = 1

$sortObj = $ws1.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws1.Range("C4:C134"))
$sortObj.SetRange($ws1.Range("A4:Q134"))
$sortObj.Header = 0
[void]$sortObj.Apply()

$ws1.Range("C4").Value = $origC4
$ws1.Range("C5").Value = $origC5
